$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels (translate to English)
$ws.Range("E1").Value = "Change"
$ws.Range("F1").Value = "Percentage Change"

# Zero out Change (E) and Percentage Change (F) values for non-trading / duplicated days
$ws.Range("F6").Value = 0
$ws.Range("F7").Value = 0
$ws.Range("E13").Value = $False
$ws.Range("F13").Value = 0
$ws.Range("E14").Value = $False
$ws.Range("F14").Value = 0
$ws.Range("E20").Value = $False
$ws.Range("F20").Value = 0
$ws.Range("E21").Value = $False
$ws.Range("F21").Value = 0
$ws.Range("E27").Value = $False
$ws.Range("F27").Value = 0
$ws.Range("E28").Value = $False
$ws.Range("F28").Value = 0
$ws.Range("E34").Value = $False
$ws.Range("F34").Value = 0
$ws.Range("E35").Value = $False
$ws.Range("F35").Value = 0
$ws.Range("F41").Value = 0
$ws.Range("F42").Value = 0
$ws.Range("E48").Value = $False
$ws.Range("F48").Value = 0
$ws.Range("E49").Value = $False
$ws.Range("F49").Value = 0
$ws.Range("E55").Value = $False
$ws.Range("F55").Value = 0
$ws.Range("E56").Value = $False
$ws.Range("F56").Value = 0
$ws.Range("E62").Value = $False
$ws.Range("F62").Value = 0
$ws.Range("E68").Value = $False
$ws.Range("F68").Value = 0
$ws.Range("E69").Value = $False
$ws.Range("F69").Value = 0
$ws.Range("E75").Value = $False
$ws.Range("F75").Value = 0
$ws.Range("E76").Value = $False
$ws.Range("F76").Value = 0
$ws.Range("E82").Value = $False
$ws.Range("F82").Value = 0
$ws.Range("E83").Value = $False
$ws.Range("F83").Value = 0
$ws.Range("E89").Value = $False
$ws.Range("F89").Value = 0
$ws.Range("E90").Value = $False
$ws.Range("F90").Value = 0
$ws.Range("F96").Value = 0
$ws.Range("F97").Value = 0
$ws.Range("E103").Value = $False
$ws.Range("F103").Value = 0
$ws.Range("E104").Value = $False
$ws.Range("F104").Value = 0
$ws.Range("F110").Value = 0
$ws.Range("F111").Value = 0
$ws.Range("F117").Value = 0
$ws.Range("F118").Value = 0
$ws.Range("F131").Value = 0
$ws.Range("F132").Value = 0
$ws.Range("F138").Value = 0
$ws.Range("F139").Value = 0
$ws.Range("E145").Value = $False
$ws.Range("F145").Value = 0
$ws.Range("E146").Value = $False
$ws.Range("F146").Value = 0
$ws.Range("F152").Value = 0
$ws.Range("F153").Value = 0
$ws.Range("F159").Value = 0
$ws.Range("F160").Value = 0
$ws.Range("E166").Value = $False
$ws.Range("F166").Value = 0
$ws.Range("E167").Value = $False
$ws.Range("F167").Value = 0
$ws.Range("E173").Value = $False
$ws.Range("F173").Value = 0
$ws.Range("E174").Value = $False
$ws.Range("F174").Value = 0
$ws.Range("F180").Value = 0
$ws.Range("F181").Value = 0
$ws.Range("E187").Value = $False
$ws.Range("F187").Value = 0
$ws.Range("E188").Value = $False
$ws.Range("F188").Value = 0
$ws.Range("F194").Value = 0
$ws.Range("F195").Value = 0
$ws.Range("E201").Value = $False
$ws.Range("F201").Value = 0
$ws.Range("E202").Value = $False
$ws.Range("F202").Value = 0
$ws.Range("E208").Value = $False
$ws.Range("F208").Value = 0
$ws.Range("E209").Value = $False
$ws.Range("F209").Value = 0
$ws.Range("E215").Value = $False
$ws.Range("F215").Value = 0
$ws.Range("E216").Value = $False
$ws.Range("F216").Value = 0
$ws.Range("F222").Value = 0
$ws.Range("F223").Value = 0
$ws.Range("F229").Value = 0
$ws.Range("F230").Value = 0
$ws.Range("F236").Value = 0
$ws.Range("F237").Value = 0
$ws.Range("F243").Value = 0
$ws.Range("F244").Value = 0
$ws.Range("E250").Value = $False
$ws.Range("F250").Value = 0
$ws.Range("E251").Value = $False
$ws.Range("F251").Value = 0
$ws.Range("F257").Value = 0
$ws.Range("F258").Value = 0
$ws.Range("F264").Value = 0
$ws.Range("F265").Value = 0
$ws.Range("F271").Value = 0
$ws.Range("F272").Value = 0
$ws.Range("F273").Value = 0
$ws.Range("E279").Value = $False
$ws.Range("F279").Value = 0
$ws.Range("E280").Value = $False
$ws.Range("F280").Value = 0
